$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (A=10): PC for Metralight system / Main Control cabinet (or in machine) / BK ; D now uses normal style; E11 empty but styled
$ws.Range("D11").Value = "BK"
$ws.Range("E11").Value = ""

# Row 12 (A=11): Metralight sensor for glue measurement / Machine / BK ; E12 = 192.168.2.60
$ws.Range("D12").Value = "BK"
$ws.Range("E12").Value = "192.168.2.60"

# Row 13 (A=12): B13 becomes "Gateway metralight EC1" ; D13 = BK ; E13 = 192.168.2.61
$ws.Range("B13").Value = "Gateway metralight EC1"
$ws.Range("D13").Value = "BK"
$ws.Range("E13").Value = "192.168.2.61"

# Row 14 (A=13): B14 becomes "Metralight debug PC" ; E14 = 192.168.2.62
$ws.Range("B14").Value = "Metralight debug PC"
$ws.Range("E14").Value = "192.168.2.62"

# Row 15 (A=14): B15 becomes "Metralight Reserve" ; E15 = 192.168.2.63
$ws.Range("B15").Value = "Metralight Reserve"
$ws.Range("E15").Value = "192.168.2.63"

# Row 16 (A=15): B16="Handheld barcode reader", C16="On machine", D16 unchanged "MAGNA", E16=192.168.2.64
$ws.Range("B16").Value = "Handheld barcode reader"
$ws.Range("C16").Value = "On machine"
$ws.Range("E16").Value = "192.168.2.64"

# Row 17 (A=16): B17="Robot - gluing", C17="Control cabinet robot 1"
$ws.Range("B17").Value = "Robot - gluing"
$ws.Range("C17").Value = "Control cabinet robot 1"

# Row 18 (A=17): B18="Robot - manipulation", C18="Control cabinet robot 2"
$ws.Range("B18").Value = "Robot - manipulation"
$ws.Range("C18").Value = "Control cabinet robot 2"

# Row 19 (A=18): C19="Gripper Robot 2"
$ws.Range("C19").Value = "Gripper Robot 2"

# Row 21 (A=20): C21="Carousel"
$ws.Range("C21").Value = "Carousel"

# Row 22 (A=21): B22="Camera 1 Primer detection", C22="Carousel"
$ws.Range("B22").Value = "Camera 1 Primer detection"
$ws.Range("C22").Value = "Carousel"

# Row 23 (A=22): B23="Camera 2 Primer detection", C23="Carousel"
$ws.Range("B23").Value = "Camera 2 Primer detection"
$ws.Range("C23").Value = "Carousel"

# Row 24 (A=23): B24="Camera 3 Primer detection", C24="Carousel"
$ws.Range("B24").Value = "Camera 3 Primer detection"
$ws.Range("C24").Value = "Carousel"

# Row 25 (A=24): B25="Camera 4 Primer detection", C25="Carousel"
$ws.Range("B25").Value = "Camera 4 Primer detection"
$ws.Range("C25").Value = "Carousel"

# Row 26 (A=25): B26="Handheld barcode reader", C26="Carousel"
$ws.Range("B26").Value = "Handheld barcode reader"
$ws.Range("C26").Value = "Carousel"

# Rows 27,28,33,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51: B becomes "Reserve"
foreach ($r in 27,28,33,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51) {
  $ws.Range("B$r").Value = "Reserve"
}

# Row 29 (A=28): C29="Rotary table"
$ws.Range("C29").Value = "Rotary table"

# Row 30 (A=29): B30="FESTO TERMINAL", C30="Rotary table"
$ws.Range("B30").Value = "FESTO TERMINAL"
$ws.Range("C30").Value = "Rotary table"

# Row 31 (A=30): B31="Camera for bar code reading 1a", C31="Rotary table"
$ws.Range("B31").Value = "Camera for bar code reading 1a"
$ws.Range("C31").Value = "Rotary table"

# Row 32 (A=31): B32="Camera for bar code reading 1b (reserve)", C32="Rotary table"
$ws.Range("B32").Value = "Camera for bar code reading 1b (reserve)"
$ws.Range("C32").Value = "Rotary table"

# Row 34 (A=33): B34="IP for programator PLC", C34="reserve"
$ws.Range("B34").Value = "IP for programator PLC"
$ws.Range("C34").Value = "reserve"

# Row 35 (A=34): B35="IP for programator ROBOTs", C35="reserve"
$ws.Range("B35").Value = "IP for programator ROBOTs"
$ws.Range("C35").Value = "reserve"

# Row 36 (A=35): B36="IP for programator Metralight", C36="reserve"
$ws.Range("B36").Value = "IP for programator Metralight"
$ws.Range("C36").Value = "reserve"

# Update selection/view to match new state
$ws.Range("E17").Select()
